{"js": "// Office.js (Word JavaScript API) script.\n// Applies the two content edits from the commit \"Partie faire un don termin\u00e9e\":\n//   1. Fix the capitalization of \"\u00e0 propos\" -> \"\u00c0 propos\" in the\n//      \"(Lien vers la page \u00e0 propos?)\" red/bold note.\n//   2. Flesh out the \"Faire un don\" section with two new explanatory\n//      paragraphs, inserted right after the existing intro paragraph and\n//      before the \"Mes dons\" heading.\n\n// ---- 1. \"\u00e0 propos\" -> \"\u00c0 propos\" -----------------------------------------\nconst oldNote = context.document.body.search(\"(Lien vers la page \u00e0 propos?)\", {\n  matchCase: true,\n  matchWildcards: false\n});\noldNote.load(\"items\");\nawait context.sync();\n\nif (oldNote.items.length > 0) {\n  // insertText(\"Replace\") keeps the run's existing formatting (bold/red).\n  oldNote.items[0].insertText(\"(Lien vers la page \u00c0 propos?)\", \"Replace\");\n  await context.sync();\n}\n\n// ---- 2. New paragraphs under \"Faire un don\" -------------------------------\nconst introResults = context.document.body.search(\n  \"C\\u2019est une fonctionnalit\u00e9 disponible \u00e0 tous les utilisateurs inscrits, puisque tout le monde est invit\u00e9 \u00e0 faire des dons aux collectivit\u00e9s.\",\n  { matchCase: true, matchWildcards: false }\n);\nintroResults.load(\"items\");\nawait context.sync();\n\nif (introResults.items.length > 0) {\n  // The search returns a Range; get its paragraph so we can insert after it.\n  const introParagraph = introResults.items[0].paragraphs.getFirst();\n  introParagraph.load(\"text\");\n  await context.sync();\n\n  // Insert in reverse order, each right \"After\" the intro paragraph, so the\n  // final order reads: intro -> \"Afin de faire un don...\" -> \"Ne vous\n  // inqui\u00e9tez pas...\".\n  introParagraph.insertParagraph(\n    \"Ne vous inqui\u00e9tez pas! Si vous faites une erreur en remplissant le formulaire, il suffit simplement d\\u2019acc\u00e9der \u00e0 la fonctionnalit\u00e9 Mes dons, disponible dans le menu sur la barre d\\u2019action.\",\n    \"After\"\n  );\n  introParagraph.insertParagraph(\n    \"Afin de faire un don aux organismes communautaires, les utilisateurs sont invit\u00e9s \u00e0 utiliser le syst\u00e8me de Gestion de produits. Il suffit simplement d\\u2019ins\u00e9rer toutes les informations demand\u00e9es dans le formulaire (Attention, tous les champs sont obligatoires) et de cliquer sur ajouter. L\\u2019application introduira donc le produit qui a \u00e9t\u00e9 ajout\u00e9 au syst\u00e8me et sera visible instantan\u00e9ment aux organismes.\",\n    \"After\"\n  );\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the two content edits from the commit \"Partie faire un don termin\u00e9e\":\n#   1. Fix the capitalization of \"\u00e0 propos\" -> \"\u00c0 propos\" in the\n#      \"(Lien vers la page \u00e0 propos?)\" red/bold note.\n#   2. Flesh out the \"Faire un don\" section with two new explanatory\n#      paragraphs, inserted right after the existing intro paragraph and\n#      before the \"Mes dons\" heading.\n\n$d = $word.ActiveDocument\n\n# ---- 1. \"\u00e0 propos\" -> \"\u00c0 propos\" -------------------------------------------\n$find1 = $d.Content\n$found1 = $find1.Find.Execute(\"(Lien vers la page \u00e0 propos?)\")\nif ($found1) {\n    # Assigning .Text on the found range replaces in place and keeps the\n    # run's existing formatting (bold/red).\n    $find1.Text = \"(Lien vers la page \u00c0 propos?)\"\n}\n\n# ---- 2. New paragraphs under \"Faire un don\" --------------------------------\n$find2 = $d.Content\n$found2 = $find2.Find.Execute(\"C\u2019est une fonctionnalit\u00e9 disponible \u00e0 tous les utilisateurs inscrits, puisque tout le monde est invit\u00e9 \u00e0 faire des dons aux collectivit\u00e9s.\")\nif ($found2) {\n    $introPara = $find2.Paragraphs(1)\n\n    # Insert an empty paragraph right after the intro paragraph and fill it\n    # with the first new explanatory paragraph.\n    $introPara.Range.InsertParagraphAfter()\n    $firstNewPara = $introPara.Next()\n    $firstNewPara.Range.Text = \"Afin de faire un don aux organismes communautaires, les utilisateurs sont invit\u00e9s \u00e0 utiliser le syst\u00e8me de Gestion de produits. Il suffit simplement d\u2019ins\u00e9rer toutes les informations demand\u00e9es dans le formulaire (Attention, tous les champs sont obligatoires) et de cliquer sur ajouter. L\u2019application introduira donc le produit qui a \u00e9t\u00e9 ajout\u00e9 au syst\u00e8me et sera visible instantan\u00e9ment aux organismes.\"\n\n    # Insert a second empty paragraph right after that one and fill it with\n    # the closing \"reassurance\" paragraph.\n    $firstNewPara = $introPara.Next()\n    $firstNewPara.Range.InsertParagraphAfter()\n    $secondNewPara = $firstNewPara.Next()\n    $secondNewPara.Range.Text = \"Ne vous inqui\u00e9tez pas! Si vous faites une erreur en remplissant le formulaire, il suffit simplement d\u2019acc\u00e9der \u00e0 la fonctionnalit\u00e9 Mes dons, disponible dans le menu sur la barre d\u2019action.\"\n}\n"}
